$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.318.40"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "1.875.38"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7131"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.91"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3105"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07751"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.94"
$ws.Range("E10").Value = "  -0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08513"
$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("D12").Value = "1.885.33"
$ws.Range("E12").Value = "  +1.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.217"
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7104"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.49"
$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").Value = "29.313.53"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008205"
$ws.Range("E17").Value = "  +5.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.008"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.17"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").Value = "2.134.43"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.814"
$ws.Range("E23").Value = "  -1.95%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1604"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.91"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.053"
$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  +1.00%  "

$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.398"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.317"
$ws.Range("E31").Value = "  +2.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.282"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.935"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7429"
$ws.Range("E36").Value = "  +2.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01866"
$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.714"
$ws.Range("E39").Value = "  +1.06%  "

$ws.Range("D40").Value = "1.182.21"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.385"
$ws.Range("E41").Value = "  +3.76%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8884"
$ws.Range("E42").Value = "  -1.75%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.97"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.33"
$ws.Range("E44").Value = "  +4.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "2.030.58"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.811"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5206"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000121"
$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.388"
$ws.Range("E50").Value = "  +0.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4311"
$ws.Range("E51").Value = "  +1.17%  "
